$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login_Check")

# Remove the extra "LOgin SUcce" values that were mistakenly placed in I3 and I4
$ws.Range("I3").ClearContents()
$ws.Range("I4").ClearContents()
